$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 210.8
$ws.Range("J2").Value = 182
$ws.Range("L2").Value = 182
$ws.Range("N2").Value = -408
$ws.Range("H19").Value = 19051.727
$ws.Range("I19").Value = 713.1667
$ws.Range("J19").Value = 41058
$ws.Range("K19").Value = 713.1667
$ws.Range("L19").Value = 41058
$ws.Range("M19").Value = -538.1667
$ws.Range("N19").Value = -41408
$ws.Range("H33").Value = 13930814
$ws.Range("I33").Value = 53755.215
$ws.Range("K33").Value = 53755.215
$ws.Range("M33").Value = -53526.215
$ws.Range("H87").Value = 1859998.8
$ws.Range("J87").Value = 1859998.8
$ws.Range("L87").Value = 1859998.8
$ws.Range("N87").Value = -1862494.8
$ws.Range("H88").Value = 989.7646999999999
$ws.Range("I88").Value = 1135.8
$ws.Range("J88").Value = 928.9167
$ws.Range("K88").Value = 1135.8
$ws.Range("L88").Value = 928.9167
$ws.Range("M88").Value = -729.8
$ws.Range("N88").Value = -1740.9167
$ws.Range("H90").Value = 1859998.8
$ws.Range("J90").Value = 1859998.8
$ws.Range("L90").Value = 5579996.4
$ws.Range("N90").Value = -5592476.4
$ws.Range("H91").Value = 989.7646999999999
$ws.Range("I91").Value = 1135.8
$ws.Range("J91").Value = 928.9167
$ws.Range("K91").Value = 1135.8
$ws.Range("L91").Value = 928.9167
$ws.Range("M91").Value = 268.2
$ws.Range("N91").Value = -3736.9167
$ws.Range("H103").Value = 249
$ws.Range("I103").Value = 249
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 747
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -161
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 5933.625
$ws.Range("I106").Value = 3791.2273
$ws.Range("K106").Value = 3791.2273
$ws.Range("M106").Value = -3160.2273
$ws.Range("H132").Value = 4339.351
$ws.Range("I132").Value = 4705.1143
$ws.Range("K132").Value = 14115.3429
$ws.Range("M132").Value = -11585.3429
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H137").Value = 3254.66
$ws.Range("I137").Value = 1917.35
$ws.Range("J137").Value = 3588.9875
$ws.Range("K137").Value = 5752.049999999999
$ws.Range("L137").Value = 10766.9625
$ws.Range("M137").Value = -3202.049999999999
$ws.Range("N137").Value = -15866.9625
$ws.Range("H141").Value = 2558.4
$ws.Range("I141").Value = 2998
$ws.Range("K141").Value = 8994
$ws.Range("M141").Value = -3814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2669.8572
$ws.Range("I2").Value = 2928
$ws.Range("K2").Value = 2928
$ws.Range("M2").Value = -2815
$ws.Range("H32").Value = 160444.78
$ws.Range("I32").Value = 162984.2
$ws.Range("K32").Value = 162984.2
$ws.Range("M32").Value = -162697.2
$ws.Range("H74").Value = 3014.516
$ws.Range("I74").Value = 1981
$ws.Range("K74").Value = 1981
$ws.Range("M74").Value = -1107
$ws.Range("H77").Value = 3014.516
$ws.Range("I77").Value = 1981
$ws.Range("K77").Value = 9905
$ws.Range("M77").Value = -5537
$ws.Range("H93").Value = 70000
$ws.Range("J93").Value = 70000
$ws.Range("L93").Value = 70000
$ws.Range("N93").Value = -74992
$ws.Range("H116").Value = 2669.8572
$ws.Range("I116").Value = 2928
$ws.Range("K116").Value = 2928
$ws.Range("M116").Value = -634
$ws.Range("H132").Value = 25643616
$ws.Range("J132").Value = 3175.3333
$ws.Range("L132").Value = 9525.999899999999
$ws.Range("N132").Value = -14585.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2669.8572
$ws.Range("I3").Value = 2928
$ws.Range("K3").Value = 2928
$ws.Range("M3").Value = -2814
$ws.Range("H22").Value = 274.75
$ws.Range("I22").Value = 274.75
$ws.Range("K22").Value = 274.75
$ws.Range("M22").Value = -101.75
$ws.Range("H134").Value = 2095.2
$ws.Range("I134").Value = 2102.487
$ws.Range("J134").Value = 2047.8334
$ws.Range("K134").Value = 6307.461
$ws.Range("L134").Value = 6143.5002
$ws.Range("M134").Value = -3772.461
$ws.Range("N134").Value = -11213.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 21.75
$ws.Range("J7").Value = 12.5
$ws.Range("L7").Value = 37.5
$ws.Range("N7").Value = -261.5
$ws.Range("H92").Value = 1957.6666
$ws.Range("J92").Value = 2999.6667
$ws.Range("L92").Value = 8999.000100000001
$ws.Range("N92").Value = -11495.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1328.625
$ws.Range("I2").Value = 902
$ws.Range("K2").Value = 902
$ws.Range("M2").Value = -789
$ws.Range("H102").Value = 1816.1428
$ws.Range("I102").Value = 623.8421
$ws.Range("K102").Value = 623.8421
$ws.Range("M102").Value = 998.1579
$ws.Range("H132").Value = 247221.14
$ws.Range("I132").Value = 359825.4
$ws.Range("J132").Value = 4688.923
$ws.Range("K132").Value = 1079476.2
$ws.Range("L132").Value = 14066.769
$ws.Range("M132").Value = -1076946.2
$ws.Range("N132").Value = -19126.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8068015.5
$ws.Range("I22").Value = 4666
$ws.Range("K22").Value = 4666
$ws.Range("M22").Value = -4371
$ws.Range("H27").Value = 8068015.5
$ws.Range("I27").Value = 4666
$ws.Range("K27").Value = 4666
$ws.Range("M27").Value = -4559
$ws.Range("H40").Value = 3197.625
$ws.Range("I40").Value = 2837.15
$ws.Range("K40").Value = 2837.15
$ws.Range("M40").Value = -2701.15
$ws.Range("H46").Value = 3062.3076
$ws.Range("I46").Value = 2796.6667
$ws.Range("J46").Value = 3142
$ws.Range("K46").Value = 2796.6667
$ws.Range("L46").Value = 3142
$ws.Range("M46").Value = -2608.6667
$ws.Range("N46").Value = -3518
$ws.Range("H68").Value = 10542.3
$ws.Range("I68").Value = 990
$ws.Range("K68").Value = 990
$ws.Range("M68").Value = -241
$ws.Range("H71").Value = 10542.3
$ws.Range("I71").Value = 990
$ws.Range("K71").Value = 4950
$ws.Range("M71").Value = -1206
$ws.Range("H82").Value = 6898
$ws.Range("I82").Value = 4697.5
$ws.Range("J82").Value = 13499.5
$ws.Range("K82").Value = 4697.5
$ws.Range("L82").Value = 13499.5
$ws.Range("M82").Value = -4336.5
$ws.Range("N82").Value = -14221.5
$ws.Range("H85").Value = 6898
$ws.Range("I85").Value = 4697.5
$ws.Range("J85").Value = 13499.5
$ws.Range("K85").Value = 4697.5
$ws.Range("L85").Value = 13499.5
$ws.Range("M85").Value = -3449.5
$ws.Range("N85").Value = -15995.5
$ws.Range("H87").Value = 73390.336
$ws.Range("I87").Value = 30171
$ws.Range("J87").Value = 95000
$ws.Range("K87").Value = 30171
$ws.Range("L87").Value = 95000
$ws.Range("M87").Value = -29048
$ws.Range("N87").Value = -97246
$ws.Range("H90").Value = 73390.336
$ws.Range("I90").Value = 30171
$ws.Range("J90").Value = 95000
$ws.Range("K90").Value = 90513
$ws.Range("L90").Value = 285000
$ws.Range("M90").Value = -84897
$ws.Range("N90").Value = -296232
$ws.Range("H122").Value = 3590.1614
$ws.Range("I122").Value = 2989.3635
$ws.Range("J122").Value = 5058.778
$ws.Range("K122").Value = 8968.0905
$ws.Range("L122").Value = 15176.334
$ws.Range("M122").Value = -6518.0905
$ws.Range("N122").Value = -20076.334
$ws.Range("H130").Value = 39997.4
$ws.Range("J130").Value = 39997.4
$ws.Range("L130").Value = 39997.4
$ws.Range("N130").Value = -50037.4
$ws.Range("H132").Value = 6383.4165
$ws.Range("I132").Value = 3693.2
$ws.Range("K132").Value = 11079.6
$ws.Range("M132").Value = -8549.599999999999
$ws.Range("H136").Value = 7484.385
$ws.Range("I136").Value = 7274.9165
$ws.Range("J136").Value = 9998
$ws.Range("K136").Value = 21824.7495
$ws.Range("L136").Value = 29994
$ws.Range("M136").Value = -19274.7495
$ws.Range("N136").Value = -35094

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16674977
$ws.Range("J81").Value = 22231194
$ws.Range("L81").Value = 44462388
$ws.Range("N81").Value = -44464510
$ws.Range("H84").Value = 16674977
$ws.Range("J84").Value = 22231194
$ws.Range("L84").Value = 222311940
$ws.Range("N84").Value = -222322548
$ws.Range("H132").Value = 468947.44
$ws.Range("I132").Value = 772166.9
$ws.Range("J132").Value = 5200.1177
$ws.Range("K132").Value = 2316500.7
$ws.Range("L132").Value = 15600.3531
$ws.Range("M132").Value = -2313970.7
$ws.Range("N132").Value = -20660.3531
